$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.409962058067322
$ws.Range("B1").Value = 2.949176073074341
$ws.Range("C1").Value = 5.319368362426758
$ws.Range("D1").Value = 2.110706806182861
$ws.Range("E1").Value = 1.183958888053894
